$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking cryptos table refresh: new Price (D) and Volume(1h) (E) readings.
# Some Price values are plain decimals that Excel would otherwise auto-convert
# to numbers; force those specific cells to keep a Text format so the value
# round-trips as a string, matching the source feed formatting.

$textPriceRows = 4,5,6,7,8,9,10,11,13,14,15,16,17,20,21,22,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,41,42,43,44,46,47,48,49,50
foreach ($r in $textPriceRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range("D2").Value = "29.326.85"
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").Value = "1.849.40"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("D4").Value = "1.013"
$ws.Range("E4").Value = "  +0.57%  "
$ws.Range("D5").Value = "244.37"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").Value = "0.6191"
$ws.Range("E6").Value = "  -1.63%  "
$ws.Range("D7").Value = "1.012"
$ws.Range("E7").Value = "  +0.65%  "
$ws.Range("D8").Value = "0.07463"
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").Value = "0.2949"
$ws.Range("E9").Value = "  +0.60%  "
$ws.Range("D10").Value = "23.06"
$ws.Range("E10").Value = "  -0.13%  "
$ws.Range("D11").Value = "0.07733"
$ws.Range("E11").Value = "  -0.01%  "
$ws.Range("D12").Value = "1.843.81"
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("D13").Value = "5.014"
$ws.Range("E13").Value = "  +0.38%  "
$ws.Range("D14").Value = "0.6751"
$ws.Range("E14").Value = "  +0.98%  "
$ws.Range("D15").Value = "83.21"
$ws.Range("E15").Value = "  +0.12%  "
$ws.Range("D16").Value = "0.000009106"
$ws.Range("E16").Value = "  -2.56%  "
$ws.Range("D17").Value = "5.906"
$ws.Range("E17").Value = "  -2.75%  "
$ws.Range("D18").Value = "29.300.72"
$ws.Range("E18").Value = "  +0.54%  "
$ws.Range("D19").Value = "2.083.50"
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("D20").Value = "238.89"
$ws.Range("E20").Value = "  +6.75%  "
$ws.Range("D21").Value = "12.68"
$ws.Range("E21").Value = "  +0.42%  "
$ws.Range("D22").Value = "1.013"
$ws.Range("E22").Value = "  +0.71%  "
$ws.Range("E23").Value = "  +0.79%  "
$ws.Range("E24").Value = "  +0.67%  "
$ws.Range("D25").Value = "160.18"
$ws.Range("E25").Value = "  -0.24%  "
$ws.Range("D26").Value = "0.1435"
$ws.Range("E26").Value = "  +2.27%  "
$ws.Range("D27").Value = "8.538"
$ws.Range("E27").Value = "  +0.35%  "
$ws.Range("D28").Value = "17.91"
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("D29").Value = "1.504"
$ws.Range("E29").Value = "  +0.17%  "
$ws.Range("D30").Value = "0.05620"
$ws.Range("E30").Value = "  +2.51%  "
$ws.Range("D31").Value = "4.161"
$ws.Range("E31").Value = "  +0.12%  "
$ws.Range("D32").Value = "4.124"
$ws.Range("E32").Value = "  +1.20%  "
$ws.Range("D33").Value = "1.221"
$ws.Range("E33").Value = "  +1.24%  "
$ws.Range("D34").Value = "1.852"
$ws.Range("E34").Value = "  -0.18%  "
$ws.Range("D35").Value = "0.7480"
$ws.Range("E35").Value = "  -0.45%  "
$ws.Range("D36").Value = "1.144"
$ws.Range("E36").Value = "  +0.69%  "
$ws.Range("D37").Value = "2.659"
$ws.Range("E37").Value = "  +1.72%  "
$ws.Range("D38").Value = "2.837"
$ws.Range("E38").Value = "  +2.98%  "
$ws.Range("D39").Value = "0.01787"
$ws.Range("E39").Value = "  -0.03%  "
$ws.Range("D40").Value = "1.216.97"
$ws.Range("E40").Value = "  -0.97%  "
$ws.Range("D41").Value = "6.496"
$ws.Range("E41").Value = "  -1.45%  "
$ws.Range("D42").Value = "0.9122"
$ws.Range("E42").Value = "  +1.79%  "
$ws.Range("D43").Value = "1.011"
$ws.Range("E43").Value = "  +0.64%  "
$ws.Range("D44").Value = "101.69"
$ws.Range("E44").Value = "  -0.39%  "
$ws.Range("D45").Value = "1.990.46"
$ws.Range("E45").Value = "  +0.25%  "
$ws.Range("D46").Value = "65.41"
$ws.Range("E46").Value = "  -0.31%  "
$ws.Range("D47").Value = "0.00000000123"
$ws.Range("E47").Value = "  +0.51%  "
$ws.Range("D48").Value = "0.5149"
$ws.Range("E48").Value = "  +0.61%  "
$ws.Range("D49").Value = "0.4060"
$ws.Range("E49").Value = "  +0.26%  "
$ws.Range("D50").Value = "9.192"
$ws.Range("E50").Value = "  +2.03%  "
$ws.Range("E51").Value = "  +0.59%  "
